# Applies the "complete section 3" report refresh:
#   - Sheet "Daily Visits": replace the two sample days with four newer days
#   - Sheet "Top Pages": re-rank the pages and add a new "/contact" row
#   - Sheet "Session Duration": swap the placeholder summary number for a
#     rendered per-session duration table, and shift the label up a row
$wb = $excel.ActiveWorkbook

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------------
# Sheet 1: "Daily Visits" - replace the two data rows with four new rows
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Daily Visits")

$ws1.Range("A2").NumberFormat = "@"
$ws1.Range("A2").Value = "2025-03-30"
$ws1.Range("B2").Value = 4

$ws1.Range("A3").NumberFormat = "@"
$ws1.Range("A3").Value = "2025-03-31"
$ws1.Range("B3").Value = 2

$ws1.Range("A4").NumberFormat = "@"
$ws1.Range("A4").Value = "2025-04-01"
$ws1.Range("B4").Value = 73

$ws1.Range("A5").NumberFormat = "@"
$ws1.Range("A5").Value = "2025-04-02"
$ws1.Range("B5").Value = 21

# Re-apply the original "row label" formatting (bold, bordered, centered) to
# every date cell so they all share a single consistent style, same as A1.
$ws1.Range("A1").Copy()
$ws1.Range("A2:A5").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet 2: "Top Pages" - reorder/update rows and add a new "/contact" row
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Top Pages")

$ws2.Range("A2").Value = "/about"
$ws2.Range("B2").Value = 26

$ws2.Range("A3").Value = "/product"
$ws2.Range("B3").Value = 23

$ws2.Range("A4").Value = "/checkout"
$ws2.Range("B4").Value = 20

$ws2.Range("A5").Value = "/contact"
$ws2.Range("B5").Value = 16

$ws2.Range("A6").Value = "/home"
$ws2.Range("B6").Value = 15

# The two new rows (5 & 6) need the same label styling as the existing ones.
$ws2.Range("A2").Copy()
$ws2.Range("A5:A6").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet 3: "Session Duration" - replace the summary value with a per-session
# duration table rendered as text, and move the label from A2/B2 to A1/A2.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Session Duration")

# Move the header label up into A1, keeping the styling the label used to
# have in A2 (bold / bordered / centered).
$ws3.Range("A1").Value = "Average Session Duration (s)"
$ws3.Range("A2").Copy()
$ws3.Range("A1").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# The old numeric placeholder in B1 is gone entirely.
$ws3.Range("B1").Clear()

# A2 becomes the rendered duration table (plain, unstyled text) and B2 (the
# old summary value) disappears.
$durationText = "              duration`nsession_id            `n5001        281.142857`n5002        319.000000`n5003        307.681818`n5004        286.952381`n5005        366.368421"
$ws3.Range("A2").ClearFormats()
$ws3.Range("A2").Value = $durationText
$ws3.Range("B2").Clear()

# Let the row shrink back to the default height instead of keeping the
# auto-expanded height computed for the multi-line text.
$ws3.Rows.Item(2).AutoFit()
